$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.537.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.841.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.625'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.92%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.70'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.31%  '
$ws.Range('E9').Value = '  +8.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0705'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.94%  '
$ws.Range('E11').Value = '  +2.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.107.83'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.842.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.33'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.675'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.47%  '
$ws.Range('E16').Value = '  +8.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.522.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0804'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '244.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.67%  '
$ws.Range('E22').Value = '  +13.86%  '
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.02'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.123'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('E29').Value = '  +26.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.326.66'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +36.92%  '
$ws.Range('E32').Value = '  +7.84%  '
$ws.Range('E33').Value = '  +7.23%  '
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('E35').Value = '  +1.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '95.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.06%  '
$ws.Range('E37').Value = '  +8.39%  '
$ws.Range('E38').Value = '  +8.03%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.11%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.348.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.32%  '
$ws.Range('E42').Value = '  +5.47%  '
$ws.Range('E43').Value = '  +7.44%  '
$ws.Range('E44').Value = '  +4.19%  '
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.28'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.72%  '
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.012.60'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.93%  '
